# Apply the cryptocurrency price/volume refresh described by the commit
# "Updated cryptos list on Sun Jul 16 15:22:20 UTC 2023 with GitHub Actions".
#
# Column D ("Price") and column E ("Volume(1h)") are refreshed for every data
# row (2-51). Rows 47/48 additionally swap their Coin/Link/Price/Volume values
# (EnergySwap <-> Aptos changed rank order).
#
# Price values are stored as literal text in the workbook (e.g. "30.405.57",
# "1.000") rather than numbers, so any new price that Excel would otherwise
# auto-convert to a number is written with a leading apostrophe (Excel's
# text-literal marker) to keep it a text value and preserve formatting such
# as trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.405.57"
$ws.Range("E2").Value = "  +0.27%  "

# Row 3
$ws.Range("D3").Value = "1.939.36"
$ws.Range("E3").Value = "  +0.15%  "

# Row 4
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'0.7684"
$ws.Range("E5").Value = "  +6.22%  "

# Row 6
$ws.Range("D6").Value = "'248.87"
$ws.Range("E6").Value = "  -0.71%  "

# Row 7
$ws.Range("D7").Value = "'0.9992"
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").Value = "'28.20"
$ws.Range("E8").Value = "  +0.99%  "

# Row 9
$ws.Range("D9").Value = "'0.3220"
$ws.Range("E9").Value = "  -2.88%  "

# Row 10
$ws.Range("D10").Value = "'0.07132"
$ws.Range("E10").Value = "  -1.74%  "

# Row 11
$ws.Range("D11").Value = "'0.7911"
$ws.Range("E11").Value = "  -2.43%  "

# Row 12
$ws.Range("D12").Value = "'0.08023"
$ws.Range("E12").Value = "  -0.79%  "

# Row 13
$ws.Range("D13").Value = "1.933.80"
$ws.Range("E13").Value = "  -0.30%  "

# Row 14
$ws.Range("D14").Value = "'5.404"
$ws.Range("E14").Value = "  -1.91%  "

# Row 15
$ws.Range("D15").Value = "'95.03"
$ws.Range("E15").Value = "  +0.26%  "

# Row 16
$ws.Range("D16").Value = "'14.68"
$ws.Range("E16").Value = "  -2.95%  "

# Row 17
$ws.Range("D17").Value = "30.405.28"
$ws.Range("E17").Value = "  +0.21%  "

# Row 18
$ws.Range("D18").Value = "'256.54"
$ws.Range("E18").Value = "  +0.98%  "

# Row 19
$ws.Range("D19").Value = "'0.000008050"
$ws.Range("E19").Value = "  -3.06%  "

# Row 20
$ws.Range("D20").Value = "'5.820"
$ws.Range("E20").Value = "  -1.34%  "

# Row 21
$ws.Range("D21").Value = "2.190.83"
$ws.Range("E21").Value = "  +0.00%  "

# Row 22
$ws.Range("D22").Value = "'0.9987"
$ws.Range("E22").Value = "  -0.16%  "

# Row 23
$ws.Range("D23").Value = "'0.9991"
$ws.Range("E23").Value = "  -0.14%  "

# Row 24
$ws.Range("D24").Value = "'6.832"
$ws.Range("E24").Value = "  -2.21%  "

# Row 25
$ws.Range("D25").Value = "'9.644"
$ws.Range("E25").Value = "  -1.30%  "

# Row 26
$ws.Range("D26").Value = "'164.80"
$ws.Range("E26").Value = "  +0.73%  "

# Row 27
$ws.Range("D27").Value = "'0.1354"
$ws.Range("E27").Value = "  +2.75%  "

# Row 28
$ws.Range("D28").Value = "'19.19"
$ws.Range("E28").Value = "  -0.61%  "

# Row 29
$ws.Range("D29").Value = "'2.319"
$ws.Range("E29").Value = "  -3.12%  "

# Row 30
$ws.Range("D30").Value = "'1.372"
$ws.Range("E30").Value = "  +1.52%  "

# Row 31
$ws.Range("D31").Value = "'1.531"
$ws.Range("E31").Value = "  -2.50%  "

# Row 32
$ws.Range("D32").Value = "'4.454"
$ws.Range("E32").Value = "  +0.26%  "

# Row 33
$ws.Range("D33").Value = "'4.170"
$ws.Range("E33").Value = "  -0.24%  "

# Row 34
$ws.Range("D34").Value = "'0.05232"
$ws.Range("E34").Value = "  +0.38%  "

# Row 35
$ws.Range("D35").Value = "'1.297"
$ws.Range("E35").Value = "  +0.55%  "

# Row 36
$ws.Range("D36").Value = "'0.7560"
$ws.Range("E36").Value = "  +0.61%  "

# Row 37
$ws.Range("E37").Value = "  +1.04%  "

# Row 38
$ws.Range("D38").Value = "'0.01978"
$ws.Range("E38").Value = "  -0.41%  "

# Row 39
$ws.Range("D39").Value = "'2.808"
$ws.Range("E39").Value = "  -0.54%  "

# Row 40
$ws.Range("D40").Value = "'78.81"
$ws.Range("E40").Value = "  -0.94%  "

# Row 41
$ws.Range("D41").Value = "'6.488"
$ws.Range("E41").Value = "  +1.89%  "

# Row 42
$ws.Range("D42").Value = "'0.4541"
$ws.Range("E42").Value = "  -0.23%  "

# Row 43
$ws.Range("D43").Value = "'1.997"
$ws.Range("E43").Value = "  -1.66%  "

# Row 44
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("D45").Value = "'0.8383"
$ws.Range("E45").Value = "  -1.00%  "

# Row 46
$ws.Range("D46").Value = "'102.34"
$ws.Range("E46").Value = "  +0.44%  "

# Row 47
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "'7.579"
$ws.Range("E47").Value = "  +1.30%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.843"
$ws.Range("E48").Value = "  +1.01%  "

# Row 49
$ws.Range("D49").Value = "'987.72"
$ws.Range("E49").Value = "  +11.82%  "

# Row 50
$ws.Range("D50").Value = "'37.56"
$ws.Range("E50").Value = "  +1.86%  "

# Row 51
$ws.Range("D51").Value = "'0.4194"
$ws.Range("E51").Value = "  -0.02%  "
